# Add two new rows of "Random" method data to the BIIB sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 ---
# Copy A2's formatting (date number format) onto A3 first so the new cell
# reuses the existing style index instead of minting a new number format.
$ws.Range("A2").Copy($ws.Range("A3"))
$ws.Range("A3").Value = 42600.835115740738
$ws.Range("B3").Value = "Random"
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 45
$ws.Range("I3").Value = 55
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 44
$ws.Range("M3").Value = 56

# --- Row 4 ---
$ws.Range("A2").Copy($ws.Range("A4"))
$ws.Range("A4").Value = 42600.88013888889
$ws.Range("B4").Value = "Random"
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 57
$ws.Range("I4").Value = 43
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 62
$ws.Range("M4").Value = 38

# Column A widened slightly (best-fit grew with the new longer date strings)
$ws.Columns("A").ColumnWidth = 14
